# Auto-generated edit script applying the cell-value changes described in the diff.
# Each row's changed cells (columns H-N) are updated per the target workbook state;
# cells that are removed entirely in the diff are cleared (ClearContents).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17: "One for the Road" (Leve Item ID 38956)
$ws.Range("H17").Value = 1040.238
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 1744.5
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 5233.5
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -5569.5

# ALC row 19: "Unbreak My Heart" (Leve Item ID 7015)
$ws.Range("H19").Value = 603
$ws.Range("I19").Value = 414.5
$ws.Range("J19").Value = 753.8
$ws.Range("K19").Value = 414.5
$ws.Range("L19").Value = 753.8
$ws.Range("M19").Value = -239.5
$ws.Range("N19").Value = -1103.8

# ALC row 38: "Just Give Him a Serum" (Leve Item ID 4599)
$ws.Range("H38").Value = 28.625
$ws.Range("I38").Value = 28.625
$ws.Range("K38").Value = 85.875
$ws.Range("M38").Value = 286.125

# ALC row 58: "A Matter of Vital Importance" (Leve Item ID 4606)
$ws.Range("H58").Value = 2300

# ALC row 87: "There Was a Late Fee" (Leve Item ID 10651)
$ws.Range("H87").Value = 99676.5
$ws.Range("J87").Value = 99676.5
$ws.Range("L87").Value = 99676.5
$ws.Range("N87").Value = -102172.5

# ALC row 90: "A Gate Arcane Is Dragon's Bane (L)" (Leve Item ID 10651)
$ws.Range("H90").Value = 99676.5
$ws.Range("J90").Value = 99676.5
$ws.Range("L90").Value = 299029.5
$ws.Range("N90").Value = -311509.5

# ALC row 106: "Making Your Mark" (Leve Item ID 19903)
$ws.Range("H106").Value = 1678.7142
$ws.Range("I106").Value = 1678.7142
$ws.Range("K106").Value = 1678.7142
$ws.Range("M106").Value = -1047.7142

# ALC row 107: "Another Man's Ink" (Leve Item ID 27766)
$ws.Range("H107").Value = 557.6667
$ws.Range("I107").Value = 524.8
$ws.Range("J107").Value = 722
$ws.Range("K107").Value = 524.8
$ws.Range("L107").Value = 722
$ws.Range("M107").Value = 1395.2
$ws.Range("N107").Value = -4562

# ALC row 137: "Cutting Edge of Culinary Quality" (Leve Item ID 44013)
$ws.Range("H137").Value = 2662.24
$ws.Range("I137").Value = 2007.2222
$ws.Range("J137").Value = 3030.6875
$ws.Range("K137").Value = 6021.6666
$ws.Range("L137").Value = 9092.0625
$ws.Range("M137").Value = -3471.6666
$ws.Range("N137").Value = -14192.0625

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2: "Ain't Got No Ingots" (Leve Item ID 27713)
$ws.Range("H2").Value = 1346.5714
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# ARM row 45: "Hollow Hallmarks" (Leve Item ID 27714)
$ws.Range("H45").Value = 2586
$ws.Range("I45").Value = 2268.5715
$ws.Range("K45").Value = 2268.5715
$ws.Range("M45").Value = -1891.5715

# ARM row 102: "Smells of Rich Tama-hagane" (Leve Item ID 19945)
$ws.Range("H102").Value = 1200
$ws.Range("I102").Value = 1200
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1200
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 422
$ws.Range("N102").ClearContents()

# ARM row 116: "No Scope" (Leve Item ID 27713)
$ws.Range("H116").Value = 1346.5714
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# ARM row 132: "Don't Bore Me, Ore Me" (Leve Item ID 43997)
$ws.Range("H132").Value = 1819.1875
$ws.Range("I132").Value = 1697
$ws.Range("K132").Value = 5091
$ws.Range("M132").Value = -2561

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3: "Hells Bells" (Leve Item ID 27713)
$ws.Range("H3").Value = 1346.5714
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# BSM row 20: "Smelt and Dealt" (Leve Item ID 14149)
$ws.Range("H20").Value = 3480.7144
$ws.Range("J20").Value = 1168.5
$ws.Range("L20").Value = 1168.5
$ws.Range("N20").Value = -1662.5

# BSM row 22: "Riveting Run" (Leve Item ID 5092)
$ws.Range("H22").Value = 333.66666
$ws.Range("I22").Value = 250.5
$ws.Range("K22").Value = 250.5
$ws.Range("M22").Value = -77.5

# BSM row 86: "Through Thick and Thin" (Leve Item ID 12526)
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# BSM row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Leve Item ID 12526)
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# BSM row 105: "Ingot to Wing It" (Leve Item ID 19947)
$ws.Range("H105").Value = 1299.5
$ws.Range("I105").Value = 1299.5
$ws.Range("K105").Value = 1299.5
$ws.Range("M105").Value = 447.5

# BSM row 140: "Ceremonial Teeth" (Leve Item ID 42471)
$ws.Range("H140").Value = 110753.11
$ws.Range("J140").Value = 110753.11
$ws.Range("L140").Value = 110753.11
$ws.Range("N140").Value = -121113.11

$ws = $wb.Worksheets.Item("CRP")
# CRP row 33: "Tools for the Tools" (Leve Item ID 1836)
$ws.Range("H33").Value = 1761.0834
$ws.Range("I33").Value = 1242.7142
$ws.Range("K33").Value = 1242.7142
$ws.Range("M33").Value = -863.7141999999999

# CRP row 134: "Wood You Be Quiet" (Leve Item ID 44020)
$ws.Range("H134").Value = 1874.1578
$ws.Range("I134").Value = 705.7857
$ws.Range("K134").Value = 2117.3571
$ws.Range("M134").Value = 417.6428999999998

$ws = $wb.Worksheets.Item("CUL")
# CUL row 80: "Saucy for a Suitor" (Leve Item ID 12890)
$ws.Range("H80").Value = 4002.2
$ws.Range("I80").Value = 3805.348
$ws.Range("J80").Value = 4649
$ws.Range("K80").Value = 11416.044
$ws.Range("L80").Value = 13947
$ws.Range("M80").Value = -10480.044
$ws.Range("N80").Value = -15819

# CUL row 83: "Saved by the Sauce (L)" (Leve Item ID 12890)
$ws.Range("H83").Value = 4002.2
$ws.Range("I83").Value = 3805.348
$ws.Range("J83").Value = 4649
$ws.Range("K83").Value = 34248.132
$ws.Range("L83").Value = 41841
$ws.Range("M83").Value = -29568.132
$ws.Range("N83").Value = -51201

# CUL row 88: "Don't Let It Fall Apart" (Leve Item ID 12851)
$ws.Range("H88").Value = 14998.5
$ws.Range("J88").Value = 14998.5
$ws.Range("L88").Value = 44995.5
$ws.Range("N88").Value = -45851.5

# CUL row 91: "Better Come Back with a Sandwich (L)" (Leve Item ID 12851)
$ws.Range("H91").Value = 14998.5
$ws.Range("J91").Value = 14998.5
$ws.Range("L91").Value = 44995.5
$ws.Range("N91").Value = -47959.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 11: "A Ringing Success" (Leve Item ID 4422)
$ws.Range("H11").Value = 9540413
$ws.Range("I11").Value = 3002512.2
$ws.Range("J11").Value = 31333416
$ws.Range("K11").Value = 3002512.2
$ws.Range("L11").Value = 31333416
$ws.Range("M11").Value = -3002373.2
$ws.Range("N11").Value = -31333694

# GSM row 13: "A Needle Is a Small Sword" (Leve Item ID 2443)
$ws.Range("H13").Value = 156.5
$ws.Range("I13").Value = 66
$ws.Range("J13").Value = 181.18182
$ws.Range("K13").Value = 66
$ws.Range("L13").Value = 181.18182
$ws.Range("M13").Value = 73
$ws.Range("N13").Value = -459.18182

# GSM row 70: "Sky Is the Limit" (Leve Item ID 14146)
$ws.Range("H70").Value = 11129.5
$ws.Range("I70").Value = 2250
$ws.Range("J70").Value = 20009
$ws.Range("K70").Value = 2250
$ws.Range("L70").Value = 20009
$ws.Range("M70").Value = -1980
$ws.Range("N70").Value = -20549

# GSM row 73: "Hulls of Broken Dreams (L)" (Leve Item ID 14146)
$ws.Range("H73").Value = 11129.5
$ws.Range("I73").Value = 2250
$ws.Range("J73").Value = 20009
$ws.Range("K73").Value = 2250
$ws.Range("L73").Value = 20009
$ws.Range("M73").Value = -1314
$ws.Range("N73").Value = -21881

# GSM row 80: "Needs More Prayerbell" (Leve Item ID 12521)
$ws.Range("H80").Value = 4950
$ws.Range("I80").Value = 4900
$ws.Range("K80").Value = 4900
$ws.Range("M80").Value = -3902

# GSM row 83: "With a Noise That Reaches Heaven (L)" (Leve Item ID 12521)
$ws.Range("H83").Value = 4950
$ws.Range("I83").Value = 4900
$ws.Range("K83").Value = 24500
$ws.Range("M83").Value = -19508

# GSM row 102: "Put the Metal to the Peddle" (Leve Item ID 36169)
$ws.Range("H102").Value = 1353.4117
$ws.Range("I102").Value = 1267.2667
$ws.Range("K102").Value = 1267.2667
$ws.Range("M102").Value = 354.7333000000001

# GSM row 113: "Copious Crystal Cannons" (Leve Item ID 27710)
$ws.Range("H113").Value = 9583.333000000001
$ws.Range("I113").Value = 7500
$ws.Range("K113").Value = 7500
$ws.Range("M113").Value = -5330

# GSM row 126: "Gold Rush Order" (Leve Item ID 36184)
$ws.Range("H126").Value = 6048.125
$ws.Range("I126").Value = 5137
$ws.Range("J126").Value = 7566.6665
$ws.Range("K126").Value = 15411
$ws.Range("L126").Value = 22699.9995
$ws.Range("M126").Value = -12941
$ws.Range("N126").Value = -27639.9995

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: "Skin off Their Backs" (Leve Item ID 5277)
$ws.Range("H22").Value = 2048.3635
$ws.Range("I22").Value = 1606.4
$ws.Range("J22").Value = 2416.6667
$ws.Range("K22").Value = 1606.4
$ws.Range("L22").Value = 2416.6667
$ws.Range("M22").Value = -1311.4
$ws.Range("N22").Value = -3006.6667

# LTW row 27: "Fire and Hide" (Leve Item ID 5277)
$ws.Range("H27").Value = 2048.3635
$ws.Range("I27").Value = 1606.4
$ws.Range("J27").Value = 2416.6667
$ws.Range("K27").Value = 1606.4
$ws.Range("L27").Value = 2416.6667
$ws.Range("M27").Value = -1499.4
$ws.Range("N27").Value = -2630.6667

# LTW row 136: "Respect for Br'aax" (Leve Item ID 44060)
$ws.Range("H136").Value = 2327.9285
$ws.Range("I136").Value = 2215.5
$ws.Range("K136").Value = 6646.5
$ws.Range("M136").Value = -4096.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 14: "Hat in Hand" (Leve Item ID 2658)
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# WVR row 133: "Begin with the Basics" (Leve Item ID 41869)
$ws.Range("H133").Value = 24999.5
$ws.Range("J133").Value = 24999.5
$ws.Range("L133").Value = 24999.5
$ws.Range("N133").Value = -35119.5

# WVR row 136: "Weaving the Envelope" (Leve Item ID 44031)
$ws.Range("H136").Value = 2073.9387
$ws.Range("I136").Value = 1477.7368
$ws.Range("K136").Value = 4433.2104
$ws.Range("M136").Value = -1883.2104

